$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E1").Value = "QUANTIDADE"

$values = @(
    5066,
    688,
    1472,
    4322,
    11323,
    24621,
    45925,
    73876,
    105679,
    135953,
    159350,
    175123,
    182628,
    179784,
    171240,
    159715,
    145280,
    130238,
    115939,
    102726,
    89889,
    79119,
    69986,
    61556,
    55091,
    48274,
    43156,
    38506,
    34873,
    31073,
    27784,
    25042,
    22536,
    19737,
    17555,
    15673,
    13189,
    11060,
    8994,
    7050,
    4849,
    3187,
    1821,
    727,
    262,
    20
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $values[$i]
}
